$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 45174 (2023-09-05)
# to 45175 (2023-09-06), keeping the existing date number format/style.
$ws.Range("C2:C5").Value = 45175
